$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.2759411334991455
$ws.Range("A3").Value = 0.07863998413085938
$ws.Range("A4").Value = 0.13780999183654785
$ws.Range("A5").Value = 0.41115498542785645
$ws.Range("A6").Value = 1.5225379467010498
$ws.Range("A7").Value = 0.2529301643371582
$ws.Range("A8").Value = 0.19219207763671875
$ws.Range("A9").Value = 0.2734339237213135
$ws.Range("A10").Value = 0.4645111560821533
$ws.Range("A11").Value = 0.20316219329833984
$ws.Range("A12").Value = 0.748046875
$ws.Range("A13").Value = 1.5321409702301025
$ws.Range("A14").Value = 0.4887099266052246
$ws.Range("A15").Value = 2.9096028804779053
$ws.Range("A16").Value = 1.4032399654388428
$ws.Range("A17").Value = 0.8683099746704102
$ws.Range("A18").Value = 6.363943815231323
$ws.Range("A19").Value = 16.865864038467407
$ws.Range("A20").Value = 1.0599842071533203
$ws.Range("A21").Value = 6.800781011581421
$ws.Range("A22").Value = 12.48893404006958
$ws.Range("A23").Value = 2.1718809604644775
$ws.Range("A24").Value = 300.0121958255768
$ws.Range("B24").Value = 0.001807459628008257
$ws.Range("C24").Value = 6350.000000000233
$ws.Range("E24").Value = 0.001807459628008257
$ws.Range("K24").Value = 6350.000000000233
$ws.Range("L24").Value = 0.03537106513977051
$ws.Range("M24").Value = 0.001807459628008257
$ws.Range("N24").Value = 6350.000000000233
$ws.Range("A25").Value = 300.0148129463196
$ws.Range("B25").Value = 0.006854978407351249
$ws.Range("C25").Value = 6211.0
$ws.Range("E25").Value = 0.006667150301761894
$ws.Range("K25").Value = 6210.000000000106
$ws.Range("L25").Value = 302.755588054657
$ws.Range("M25").Value = 0.006667150301761894
$ws.Range("N25").Value = 6210.000000000106
$ws.Range("O25").Value = 0.11976122856140137
$ws.Range("P25").Value = 0.006667150301761894
$ws.Range("Q25").Value = 6210.000000000106
$ws.Range("A26").Value = 19.12791895866394
$ws.Range("A27").Value = 300.0176281929016
$ws.Range("B27").Value = 0.0011054114588942687
$ws.Range("C27").Value = 12438.00000000025
$ws.Range("E27").Value = 0.0011054114588942687
$ws.Range("K27").Value = 12438.00000000025
$ws.Range("L27").Value = 0.11071395874023438
$ws.Range("M27").Value = 0.0011054114588942687
$ws.Range("N27").Value = 12438.00000000025
$ws.Range("A28").Value = 303.4111669063568
$ws.Range("B28").Value = 0.004145625336440365
$ws.Range("E28").Value = 0.003739271824528996
$ws.Range("K28").Value = 12272.0
$ws.Range("L28").Value = 46.88124918937683
$ws.Range("M28").Value = 0.003739271824528996
$ws.Range("N28").Value = 12272.0
$ws.Range("O28").Value = $null
$ws.Range("P28").Value = $null
$ws.Range("Q28").Value = $null
$ws.Range("A29").Value = 22.52096199989319
$ws.Range("A30").Value = 300.032438993454
$ws.Range("B30").Value = 0.0015431161130305122
$ws.Range("E30").Value = 0.0015431161130305122
$ws.Range("L30").Value = 0.07541108131408691
$ws.Range("M30").Value = 0.0015431161130305122
$ws.Range("A31").Value = 301.4178650379181
$ws.Range("B31").Value = 0.0052293990276651055
$ws.Range("C31").Value = 24477.000000000146
$ws.Range("E31").Value = 0.00510746097900461
$ws.Range("K31").Value = 24474.000000000146
$ws.Range("L31").Value = 301.9986660480499
$ws.Range("M31").Value = 0.00510746097900461
$ws.Range("N31").Value = 24474.000000000146
$ws.Range("O31").Value = 0.057090044021606445
$ws.Range("P31").Value = 0.00510746097900461
$ws.Range("Q31").Value = 24474.000000000146
$ws.Range("A32").Value = 17.767588138580322
$ws.Range("A33").Value = 304.35000801086426
$ws.Range("B33").Value = 0.0012933204289859999
$ws.Range("E33").Value = 0.0012933204289859999
$ws.Range("L33").Value = 0.04673504829406738
$ws.Range("M33").Value = 0.0012933204289859999
$ws.Range("A34").Value = 300.22461104393005
$ws.Range("B34").Value = 0.003907605222313116
$ws.Range("C34").Value = 54765.0
$ws.Range("E34").Value = 0.003907605222313116
$ws.Range("K34").Value = 54765.0
$ws.Range("L34").Value = 0.03720402717590332
$ws.Range("M34").Value = 0.003907605222313116
$ws.Range("N34").Value = 54765.0
$ws.Range("A35").Value = 151.76260113716125
$ws.Range("A36").Value = 300.25423312187195
$ws.Range("L36").Value = 0.03326892852783203
$ws.Range("A37").Value = 300.3874320983887
$ws.Range("L37").Value = 0.8923239707946777
$ws.Range("A38").Value = 0.7792611122131348
$ws.Range("A39").Value = 2.793440103530884
$ws.Range("A40").Value = 5.638268232345581
$ws.Range("A41").Value = 0.43352794647216797
$ws.Range("A42").Value = 1.0116770267486572
$ws.Range("A43").Value = 5.551558017730713
$ws.Range("A44").Value = 2.148998975753784
$ws.Range("A45").Value = 5.780341863632202
$ws.Range("A46").Value = 43.02309012413025
$ws.Range("A47").Value = 1.7529680728912354
$ws.Range("A48").Value = 3.6590211391448975
$ws.Range("A49").Value = 130.4853971004486
$ws.Range("A50").Value = 1.848512887954712
$ws.Range("A51").Value = 10.781486988067627
$ws.Range("A52").Value = 218.38122391700745
